$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$aValues = @(5660,5610,5580,5550,5520,5500,5490,5470,5450,5430,5430,5440,5450,5460,5470,5500,5560,5630,5720,5820,5940,6060,6190,6360,6540,6730,6890,7020,7130,7180,7190,7200,7180,7150,7090,7000,6900,6800,6700,6600,6480,6390,6310,6240,6170,6120,6080,6050,6030,6030,6020,6020,6020,6020,6030,6030,6040,6050,6060,6070,6110,6150,6210,6270,6350,6420,6490,6570,6650,6730,6820,6930,7040,7160,7280,7400,7540,7610,7620,7620,7570,7480,7380,7250,7050,6900,6720,6550,6400,6260,6140,6020,5930,5860,5810,5750)
$bValues = @(45751,45751.01041666666,45751.02083333334,45751.03125,45751.04166666666,45751.05208333334,45751.0625,45751.07291666666,45751.08333333334,45751.09375,45751.10416666666,45751.11458333334,45751.125,45751.13541666666,45751.14583333334,45751.15625,45751.16666666666,45751.17708333334,45751.1875,45751.19791666666,45751.20833333334,45751.21875,45751.22916666666,45751.23958333334,45751.25,45751.26041666666,45751.27083333334,45751.28125,45751.29166666666,45751.30208333334,45751.3125,45751.32291666666,45751.33333333334,45751.34375,45751.35416666666,45751.36458333334,45751.375,45751.38541666666,45751.39583333334,45751.40625,45751.41666666666,45751.42708333334,45751.4375,45751.44791666666,45751.45833333334,45751.46875,45751.47916666666,45751.48958333334,45751.5,45751.51041666666,45751.52083333334,45751.53125,45751.54166666666,45751.55208333334,45751.5625,45751.57291666666,45751.58333333334,45751.59375,45751.60416666666,45751.61458333334,45751.625,45751.63541666666,45751.64583333334,45751.65625,45751.66666666666,45751.67708333334,45751.6875,45751.69791666666,45751.70833333334,45751.71875,45751.72916666666,45751.73958333334,45751.75,45751.76041666666,45751.77083333334,45751.78125,45751.79166666666,45751.80208333334,45751.8125,45751.82291666666,45751.83333333334,45751.84375,45751.85416666666,45751.86458333334,45751.875,45751.88541666666,45751.89583333334,45751.90625,45751.91666666666,45751.92708333334,45751.9375,45751.94791666666,45751.95833333334,45751.96875,45751.97916666666,45751.98958333334)

for ($i = 0; $i -lt $aValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $aValues[$i]
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
}
